$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Fill in login form data
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Zoom in on the sheet
$excel.ActiveWindow.Zoom = 205

# Set selection to B3 (where the cursor ends up after data entry)
$ws.Range("B3").Select()
